# gestionActivos.xlsx update
# - Re-organizes / enumerates the test case steps for CP_GESACT_003 and merges
#   the former CP_GESACT_004 "LOST update" case content into it, adding the
#   FACTORYSERIAL filter verification steps.
# - Updates row 4 description / expected / obtained result texts accordingly.
# - Adjusts row 4 height and the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (CP_GESACT_003) cell content updates ---

# F4: Pasos a seguir (renumbered / expanded steps)
$pasos = @'
1.Seleccionar el primer registro de la tabla
2.Clic en el botón "Actualizar estado operativo"
3.Clic en estado para mostrar la lista opciones.
4.Seleccionar la opción "LOST" en el estado
5.Diligenciar el campo de comentario con "test automatización"
6.Clic en botón "Guardar" para actualizar estado ont
7.Clic en botón "Seleccionar entidad"
8.Seleccionar la fila que contiene el texto "elemento secundario"
9.Clic en botón "Siguiente"
10.Seleccionar "ont"
11.Clic en botón "Siguiente"
12.Seleccionar fila con ID 10 "LOST"
13.Hacer clic en el botón "FINALIZAR"
14.Clic en el botón "Mostrar filtro"
15.Clic en el estado para mostrar la lista de opciones.
16.Seleccionar opción "FACTORYSERIAL" en la lista de opciones
17.Diligenciar campo de texto con serial ONT capturado previamente
18.Clic en botón "Aplicar filtro"
'@
$ws.Range("F4").Value = $pasos

# G4 (Datos de prueba) stays "Serial ONT válido" - no change needed

# H4: Resultado esperado
$ws.Range("H4").Value = 'El estado FAILED del dispositivo se actualiza a LOST y se cierra el modal.'

# I4: Resultado obtenido
$ws.Range("I4").Value = 'La tabla se filtró correctamente mostrando solo el registro con el FACTORYSERIAL indicado.'

# B4: Nombre/Descripcion
$ws.Range("B4").Value = 'Actualiza ONT a LOST y valida filtrado por FACTORYSERIAL'

# --- Row height adjustment for row 4 ---
$ws.Rows.Item(4).RowHeight = 124.5

# --- Update active selection to C8 (matches saved sheetView selection) ---
$ws.Range("C8").Select()
